$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-6 from 2023-09-01 (45170)
# to 2023-09-05 (45174), keeping existing formatting.
foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = 45174
}
